$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 173.8357176670649
$ws.Range("C2").Value = 223.70860614204946
$ws.Range("D2").Value = 175.68100970774961
$ws.Range("E2").Value = 225.06406872434579

# Row 3 values
$ws.Range("B3").Value = 171.82536395782364
$ws.Range("C3").Value = 229.76529476976782
$ws.Range("D3").Value = 183.72897075034024
$ws.Range("E3").Value = 219.50334509534949

# Update selection to match new selection range
$ws.Range("B1:E3").Select()
